$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-02 Friday" "2024-08-03 Saturday"

Replace-Text "27×65=" "69×26="
Replace-Text "98×11=" "89×33="
Replace-Text "37×47=" "47×27="
Replace-Text "88×21=" "38×35="
Replace-Text "91×27=" "60×69="
Replace-Text "40×30=" "77×19="
Replace-Text "79×63=" "33×87="
Replace-Text "84×97=" "72×28="
Replace-Text "57×28=" "40×71="
Replace-Text "94×87=" "64×99="
Replace-Text "65×18=" "17×75="
Replace-Text "20×72=" "78×44="
Replace-Text "73×88=" "64×81="
Replace-Text "59×93=" "93×27="
Replace-Text "99×25=" "74×80="
Replace-Text "78×47=" "84×24="
Replace-Text "30×98=" "87×22="
Replace-Text "41×76=" "52×71="
Replace-Text "33×62=" "34×29="
Replace-Text "32×90=" "85×90="
Replace-Text "65×41=" "94×48="
Replace-Text "89×86=" "72×14="
Replace-Text "71×25=" "70×51="
Replace-Text "65×94=" "65×81="
Replace-Text "94×85=" "38×51="
